$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 17 (keep only header row 1 and data row 2)
$ws.Range("A3:B17").EntireRow.Delete()

# Update the remaining data row values
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 0.8557439673732903
